$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.860.73"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.629.37"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "1.627.05"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "27.870.74"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.76%  "
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.417.01"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E35").Value = "  +2.56%  "
$ws.Range("E36").Value = "  -5.42%  "
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0169"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.553"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("D45").Value = "1.770.50"
$ws.Range("E46").Value = "  -4.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("E51").Value = "  -0.43%  "
